$d = $word.ActiveDocument

# Find the bullet "Jumlah total karakter (tanpa spasi)." and, inside
# it, the word "tanpa" that was split by the author's in-place edit.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jumlah total karakter*tanpa*") {
        $target = $p.Range
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Jumlah total karakter (tanpa spasi)' bullet"
}

$found = $target.Find.Execute("tanpa", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'tanpa' in the target bullet"
}

# $target now covers exactly "tanpa" (Find collapses/extends it to the
# match). Split it after "ta" -> "ta" | "npa".
$splitAt = $target.Start + 2
$splitPoint = $d.Range($splitAt, $splitAt)

# Re-dropping the _GoBack bookmark here both creates the split (Word
# always keeps the bookmark's host run boundary at the bookmark) and
# moves it away from its old home at the top of the document, since a
# document can only have one bookmark with a given name.
$d.Bookmarks.Add("_GoBack", $splitPoint)
$bm = $d.Bookmarks("_GoBack")

# Rewrite the trailing "npa" so Word treats it as freshly authored text
# (no rsid) rather than merging it back with the untouched "ta" run -
# matching the author's actual edit-in-place-then-retype sequence.
$tail = $d.Range($bm.End, $bm.End + 3)
$tail.Text = "zzz"
$tail = $d.Range($bm.End, $bm.End + 3)
$tail.Text = "npa"
